# Update "想去人数" (interest count) figures in column F across the
# 展览 / 演出 / 全部类型 sheets to match the freshly re-scraped data.
$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 916
$ws1.Range("F3").Value = 512
$ws1.Range("F4").Value = 512
$ws1.Range("F5").Value = 799
$ws1.Range("F6").Value = 191
$ws1.Range("F7").Value = 1386
$ws1.Range("F8").Value = 820
$ws1.Range("F9").Value = 463
$ws1.Range("F10").Value = 622
$ws1.Range("F11").Value = 177
$ws1.Range("F13").Value = 32
$ws1.Range("F14").Value = 206
$ws1.Range("F15").Value = 107
$ws1.Range("F16").Value = 1559
$ws1.Range("F17").Value = 195
$ws1.Range("F18").Value = 28
$ws1.Range("F19").Value = 465
$ws1.Range("F20").Value = 61
$ws1.Range("F21").Value = 390
$ws1.Range("F23").Value = 629
$ws1.Range("F24").Value = 23
$ws1.Range("F25").Value = 211
$ws1.Range("F26").Value = 718
$ws1.Range("F28").Value = 1398
$ws1.Range("F29").Value = 147

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 31
$ws2.Range("F6").Value = 8
$ws2.Range("F9").Value = 36

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 916
$ws4.Range("F4").Value = 512
$ws4.Range("F5").Value = 512
$ws4.Range("F6").Value = 799
$ws4.Range("F7").Value = 191
$ws4.Range("F8").Value = 1386
$ws4.Range("F9").Value = 820
$ws4.Range("F11").Value = 31
$ws4.Range("F12").Value = 463
$ws4.Range("F13").Value = 622
$ws4.Range("F15").Value = 177
$ws4.Range("F17").Value = 32
$ws4.Range("F18").Value = 206
$ws4.Range("F19").Value = 107
$ws4.Range("F20").Value = 1559
$ws4.Range("F22").Value = 195
$ws4.Range("F23").Value = 28
$ws4.Range("F24").Value = 465
$ws4.Range("F25").Value = 61
$ws4.Range("F26").Value = 390
$ws4.Range("F27").Value = 8
$ws4.Range("F31").Value = 629
$ws4.Range("F32").Value = 36
$ws4.Range("F36").Value = 23
$ws4.Range("F37").Value = 211
$ws4.Range("F38").Value = 718
$ws4.Range("F40").Value = 1398
$ws4.Range("F41").Value = 147
